$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.6151236666666667
$ws.Range("H2").Value2 = 1.845371
$ws.Range("I2").Value2 = 0.01505763605988265
$ws.Range("J2").Value2 = 0.01505763605988265
$ws.Range("M2").Value2 = 19.21315233333334
$ws.Range("N2").Value2 = 57.63945700000001
$ws.Range("O2").Value2 = 0.04451179209991234
$ws.Range("P2").Value2 = 0.04451179209991233
$ws.Range("Q2").Value2 = 11.81846471150523
$ws.Range("R2").Value2 = 106.366182403547
$ws.Range("S2").Value2 = 0.0006702423658136397
$ws.Range("T2").Value2 = 0.0006702423658136396
# Row 3
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.6151236666666667
$ws.Range("H3").Value2 = 1.845371
$ws.Range("I3").Value2 = 0.01505763605988265
$ws.Range("J3").Value2 = 0.01505763605988265
$ws.Range("O3").Value2 = 0.2141755495962477
$ws.Range("P3").Value2 = 0.2141755495962477
$ws.Range("Q3").Value2 = 56.86641798849245
$ws.Range("R3").Value2 = 511.797761896432
$ws.Range("S3").Value2 = 0.003224977478745645
$ws.Range("T3").Value2 = 0.003224977478745645
# Row 4
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.6151236666666667
$ws.Range("H4").Value2 = 1.845371
$ws.Range("I4").Value2 = 0.01505763605988265
$ws.Range("J4").Value2 = 0.01505763605988265
$ws.Range("M4").Value2 = 166.8580016666666
$ws.Range("N4").Value2 = 500.5740049999999
$ws.Range("O4").Value2 = 0.3865658561145097
$ws.Range("P4").Value2 = 0.3865658561145097
$ws.Range("Q4").Value2 = 102.6383057978728
$ws.Range("R4").Value2 = 923.7447521808549
$ws.Range("S4").Value2 = 0.00582076797454925
$ws.Range("T4").Value2 = 0.00582076797454925
# Row 5
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.6151236666666667
$ws.Range("H5").Value2 = 1.845371
$ws.Range("I5").Value2 = 0.01505763605988265
$ws.Range("J5").Value2 = 0.01505763605988265
$ws.Range("M5").Value2 = 41.09915599999999
$ws.Range("N5").Value2 = 123.297468
$ws.Range("O5").Value2 = 0.09521587377309249
$ws.Range("P5").Value2 = 0.09521587377309249
$ws.Range("Q5").Value2 = 25.28106353562533
$ws.Range("R5").Value2 = 227.529571820628
$ws.Range("S5").Value2 = 0.001433725974398952
$ws.Range("T5").Value2 = 0.001433725974398952
# Row 6
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.6151236666666667
$ws.Range("H6").Value2 = 1.845371
$ws.Range("I6").Value2 = 0.01505763605988265
$ws.Range("J6").Value2 = 0.01505763605988265
$ws.Range("M6").Value2 = 112.0244103333333
$ws.Range("N6").Value2 = 336.073231
$ws.Range("O6").Value2 = 0.2595309284162377
$ws.Range("P6").Value2 = 0.2595309284162377
$ws.Range("Q6").Value2 = 68.90886604041123
$ws.Range("R6").Value2 = 620.179794363701
$ws.Range("S6").Value2 = 0.003907922266375164
$ws.Range("T6").Value2 = 0.003907922266375164
# Row 7
$ws.Range("I7").Value2 = 0.9124571722898065
$ws.Range("J7").Value2 = 0.9124571722898065
$ws.Range("M7").Value2 = 19.21315233333334
$ws.Range("N7").Value2 = 57.63945700000001
$ws.Range("O7").Value2 = 0.04451179209991234
$ws.Range("P7").Value2 = 0.04451179209991233
$ws.Range("Q7").Value2 = 716.1710409642457
$ws.Range("R7").Value2 = 6445.539368678211
$ws.Range("S7").Value2 = 0.04061510395303776
$ws.Range("T7").Value2 = 0.04061510395303775
# Row 8
$ws.Range("I8").Value2 = 0.9124571722898065
$ws.Range("J8").Value2 = 0.9124571722898065
$ws.Range("O8").Value2 = 0.2141755495962477
$ws.Range("P8").Value2 = 0.2141755495962477
$ws.Range("S8").Value2 = 0.1954260163582074
$ws.Range("T8").Value2 = 0.1954260163582074
# Row 9
$ws.Range("I9").Value2 = 0.9124571722898065
$ws.Range("J9").Value2 = 0.9124571722898065
$ws.Range("M9").Value2 = 166.8580016666666
$ws.Range("N9").Value2 = 500.5740049999999
$ws.Range("O9").Value2 = 0.3865658561145097
$ws.Range("P9").Value2 = 0.3865658561145097
$ws.Range("Q9").Value2 = 6219.638853303067
$ws.Range("R9").Value2 = 55976.7496797276
$ws.Range("S9").Value2 = 0.3527247879740337
$ws.Range("T9").Value2 = 0.3527247879740337
# Row 10
$ws.Range("I10").Value2 = 0.9124571722898065
$ws.Range("J10").Value2 = 0.9124571722898065
$ws.Range("M10").Value2 = 41.09915599999999
$ws.Range("N10").Value2 = 123.297468
$ws.Range("O10").Value2 = 0.09521587377309249
$ws.Range("P10").Value2 = 0.09521587377309249
$ws.Range("Q10").Value2 = 1531.972724965396
$ws.Range("R10").Value2 = 13787.75452468856
$ws.Range("S10").Value2 = 0.08688040694009912
$ws.Range("T10").Value2 = 0.08688040694009912
# Row 11
$ws.Range("I11").Value2 = 0.9124571722898065
$ws.Range("J11").Value2 = 0.9124571722898065
$ws.Range("M11").Value2 = 112.0244103333333
$ws.Range("N11").Value2 = 336.073231
$ws.Range("O11").Value2 = 0.2595309284162377
$ws.Range("P11").Value2 = 0.2595309284162377
$ws.Range("Q11").Value2 = 4175.714488175823
$ws.Range("R11").Value2 = 37581.43039358241
$ws.Range("S11").Value2 = 0.2368108570644284
$ws.Range("T11").Value2 = 0.2368108570644284
# Row 12
$ws.Range("G12").Value2 = 2.961112666666666
$ws.Range("H12").Value2 = 8.883337999999998
$ws.Range("I12").Value2 = 0.07248519165031087
$ws.Range("J12").Value2 = 0.07248519165031085
$ws.Range("M12").Value2 = 19.21315233333334
$ws.Range("N12").Value2 = 57.63945700000001
$ws.Range("O12").Value2 = 0.04451179209991234
$ws.Range("P12").Value2 = 0.04451179209991233
$ws.Range("Q12").Value2 = 56.89230874082956
$ws.Range("R12").Value2 = 512.0307786674659
$ws.Range("S12").Value2 = 0.003226445781060939
$ws.Range("T12").Value2 = 0.003226445781060938
# Row 13
$ws.Range("G13").Value2 = 2.961112666666666
$ws.Range("H13").Value2 = 8.883337999999998
$ws.Range("I13").Value2 = 0.07248519165031087
$ws.Range("J13").Value2 = 0.07248519165031085
$ws.Range("O13").Value2 = 0.2141755495962477
$ws.Range("P13").Value2 = 0.2141755495962477
$ws.Range("Q13").Value2 = 273.7463696140551
$ws.Range("R13").Value2 = 2463.717326526496
$ws.Range("S13").Value2 = 0.01552455575929468
$ws.Range("T13").Value2 = 0.01552455575929467
# Row 14
$ws.Range("G14").Value2 = 2.961112666666666
$ws.Range("H14").Value2 = 8.883337999999998
$ws.Range("I14").Value2 = 0.07248519165031087
$ws.Range("J14").Value2 = 0.07248519165031085
$ws.Range("M14").Value2 = 166.8580016666666
$ws.Range("N14").Value2 = 500.5740049999999
$ws.Range("O14").Value2 = 0.3865658561145097
$ws.Range("P14").Value2 = 0.3865658561145097
$ws.Range("Q14").Value2 = 494.0853422698543
$ws.Range("R14").Value2 = 4446.768080428688
$ws.Range("S14").Value2 = 0.02802030016592673
$ws.Range("T14").Value2 = 0.02802030016592673
# Row 15
$ws.Range("G15").Value2 = 2.961112666666666
$ws.Range("H15").Value2 = 8.883337999999998
$ws.Range("I15").Value2 = 0.07248519165031087
$ws.Range("J15").Value2 = 0.07248519165031085
$ws.Range("M15").Value2 = 41.09915599999999
$ws.Range("N15").Value2 = 123.297468
$ws.Range("O15").Value2 = 0.09521587377309249
$ws.Range("P15").Value2 = 0.09521587377309249
$ws.Range("Q15").Value2 = 121.6992314209093
$ws.Range("R15").Value2 = 1095.293082788184
$ws.Range("S15").Value2 = 0.006901740858594417
$ws.Range("T15").Value2 = 0.006901740858594415
# Row 16
$ws.Range("G16").Value2 = 2.961112666666666
$ws.Range("H16").Value2 = 8.883337999999998
$ws.Range("I16").Value2 = 0.07248519165031087
$ws.Range("J16").Value2 = 0.07248519165031085
$ws.Range("M16").Value2 = 112.0244103333333
$ws.Range("N16").Value2 = 336.073231
$ws.Range("O16").Value2 = 0.2595309284162377
$ws.Range("P16").Value2 = 0.2595309284162377
$ws.Range("Q16").Value2 = 331.7169004138975
$ws.Range("R16").Value2 = 2985.452103725077
$ws.Range("S16").Value2 = 0.0188121490854341
$ws.Range("T16").Value2 = 0.0188121490854341
